$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 33, pushing the
# existing rows 33:51 down to 34:52.
$ws.Rows(33).Insert()

$ws.Cells.Item(33, 1).Value2 = 10
$ws.Cells.Item(33, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(33, 3).Value2 = "La Araucanía"
$ws.Cells.Item(33, 4).Value2 = 44879
$ws.Cells.Item(33, 5).Value2 = 9
$ws.Cells.Item(33, 6).Value2 = 100112042
$ws.Cells.Item(33, 7).Value2 = "Locoto"
$ws.Cells.Item(33, 8).Value2 = "Sin especificar"
$ws.Cells.Item(33, 9).Value2 = "Primera"
$ws.Cells.Item(33, 10).Value2 = 200
$ws.Cells.Item(33, 11).Value2 = 2500
$ws.Cells.Item(33, 12).Value2 = 2500
$ws.Cells.Item(33, 13).Value2 = 2500
$ws.Cells.Item(33, 14).Value2 = "$/kilo"
$ws.Cells.Item(33, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value2 = 2500
$ws.Cells.Item(33, 17).Value2 = 1
$ws.Cells.Item(33, 18).Value2 = "Hortaliza"
